$wb = $excel.ActiveWorkbook

# --- Sheet: "Estadisticos 2P" ---
$ws2 = $wb.Worksheets.Item("Estadisticos 2P")

$ws2.Range("D2").Value = 0
$ws2.Range("E2").Value = 3
$ws2.Range("F2").Value = 16
$ws2.Range("G2").Value = 84.20999999999999
$ws2.Range("H2").Value = 8.300000000000001

$ws2.Range("D3").Value = 0
$ws2.Range("E3").Value = 0
$ws2.Range("F3").Value = 25
$ws2.Range("G3").Value = 100
$ws2.Range("H3").Value = 8.199999999999999

$ws2.Range("D4").Value = 0
$ws2.Range("E4").Value = 0
$ws2.Range("F4").Value = 13
$ws2.Range("G4").Value = 100
$ws2.Range("H4").Value = 8.199999999999999

$ws2.Range("D5").Value = 0
$ws2.Range("E5").Value = 0
$ws2.Range("F5").Value = 14
$ws2.Range("G5").Value = 100
$ws2.Range("H5").Value = 9.1

$ws2.Range("D6").Value = 0
$ws2.Range("E6").Value = 0
$ws2.Range("F6").Value = 10
$ws2.Range("G6").Value = 100
$ws2.Range("H6").Value = 9.4

# --- Sheet: "Estadisticos Final" ---
$ws3 = $wb.Worksheets.Item("Estadisticos Final")

$ws3.Range("E2").Value = 3
$ws3.Range("F2").Value = 16
$ws3.Range("G2").Value = 84.20999999999999
$ws3.Range("H2").Value = 7.9

$ws3.Range("E3").Value = 0
$ws3.Range("F3").Value = 25
$ws3.Range("G3").Value = 100
$ws3.Range("H3").Value = 8.300000000000001

$ws3.Range("E4").Value = 0
$ws3.Range("F4").Value = 13
$ws3.Range("G4").Value = 100

$ws3.Range("H6").Value = 9.300000000000001

# --- Sheet: "Rescatables" ---
$ws4 = $wb.Worksheets.Item("Rescatables")
$ws4.Rows.Item(2).Delete()
